$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: payment method corrected to CASH, registered balance updated ---
$ws.Range("J17").Value = "CASH"
$ws.Range("N17").Value = 2324735

# --- New rows 18-25 ---
# Row 18
$ws.Range("A18").Value = "64bba574a7499d7733f537e3"
$ws.Range("B18").Value = "640900b19a139999b1824b31"
$ws.Range("C18").Value = 4200673
$ws.Range("D18").Value = 1900
$ws.Range("E18").Value = 45129.04207175926
$ws.Range("F18").Value = "Cash payment"
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = "Invoiced"
$ws.Range("J18").Value = "CASH"
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M18").Value = "Bernard"
$ws.Range("N18").Value = 4200673
$ws.Range("O18").Value = "'976105"
$ws.Range("P18").Value = "Approved"

# Row 19
$ws.Range("A19").Value = "64bc3de5863ca7a448d1955e"
$ws.Range("B19").Value = "6405f77c1ec640bb7919b533"
$ws.Range("C19").Value = 1884668
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 45129.04207175926
$ws.Range("F19").Value = "Cash payment"
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = "Invoiced"
$ws.Range("J19").Value = "CASH"
$ws.Range("K19").Value = 10
$ws.Range("L19").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M19").Value = "Bernard"
$ws.Range("N19").Value = 1884668
$ws.Range("O19").Value = "'566965"
$ws.Range("P19").Value = "Approved"

# Row 20
$ws.Range("A20").Value = "64bc62d7250f0d99f9814523"
$ws.Range("B20").Value = "640900b19a139999b1824b31"
$ws.Range("C20").Value = 1422049
$ws.Range("D20").Value = 1200
$ws.Range("E20").Value = 45131.04207175926
$ws.Range("F20").Value = "Bank Payment"
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = "Invoiced"
$ws.Range("J20").Value = "GT Bank"
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = "64996a6f08c70837359160b9"
$ws.Range("M20").Value = "Bernard"
$ws.Range("N20").Value = 1422049
$ws.Range("O20").Value = "'519375"
$ws.Range("P20").Value = "Approved"

# Row 21
$ws.Range("A21").Value = "64bdb5afa815a7d13319bc28"
$ws.Range("B21").Value = "640900b19a139999b1824b31"
$ws.Range("C21").Value = 20260
$ws.Range("D21").Value = 1000
$ws.Range("E21").Value = 45131.04207175926
$ws.Range("F21").Value = "Bank Payment"
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = "Invoiced"
$ws.Range("J21").Value = "GT Bank"
$ws.Range("K21").Value = 10
$ws.Range("L21").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M21").Value = "Bernard"
$ws.Range("N21").Value = 20260
$ws.Range("O21").Value = "'98492"
$ws.Range("P21").Value = "Approved"

# Row 22
$ws.Range("A22").Value = "64bdd2d5344b2b124336df12"
$ws.Range("B22").Value = "640900b19a139999b1824b31"
$ws.Range("C22").Value = 930494
$ws.Range("D22").Value = 1000
$ws.Range("E22").Value = 45131.04207175926
$ws.Range("F22").Value = "Cash payment"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = "Invoiced"
$ws.Range("J22").Value = "CASH"
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M22").Value = "Bernard"
$ws.Range("N22").Value = 930494
$ws.Range("O22").Value = "'325346"
$ws.Range("P22").Value = "Approved"

# Row 23
$ws.Range("A23").Value = "64be3289c3b8ea80b088774f"
$ws.Range("B23").Value = "640900b19a139999b1824b31"
$ws.Range("C23").Value = 1395241
$ws.Range("D23").Value = 1000
$ws.Range("E23").Value = 45131.04207175926
$ws.Range("F23").Value = "Bank Payment"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = "Invoiced"
$ws.Range("J23").Value = "GT Bank"
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M23").Value = "Bernard"
$ws.Range("N23").Value = 1395241
$ws.Range("O23").Value = "'161058"
$ws.Range("P23").Value = "Approved"

# Row 24
$ws.Range("A24").Value = "64c03a836452e19bacfb6b82"
$ws.Range("B24").Value = "640900b19a139999b1824b31"
$ws.Range("C24").Value = 28680
$ws.Range("D24").Value = 3000
$ws.Range("E24").Value = 45132.04207175926
$ws.Range("F24").Value = "Cash payment"
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = "Invoiced"
$ws.Range("J24").Value = "CASH"
$ws.Range("K24").Value = 10
$ws.Range("L24").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M24").Value = "Bernard"
$ws.Range("N24").Value = 28680
$ws.Range("O24").Value = "'989329"
$ws.Range("P24").Value = "Approved"

# Row 25
$ws.Range("A25").Value = "64c21fb6bd52a08c991a112a"
$ws.Range("B25").Value = "640900b19a139999b1824b31"
$ws.Range("C25").Value = 4209403
$ws.Range("D25").Value = 1000
$ws.Range("E25").Value = 45138.04207175926
$ws.Range("F25").Value = "Cash payment"
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = "Invoiced"
$ws.Range("J25").Value = "GT Bank"
$ws.Range("K25").Value = 10
$ws.Range("L25").Value = "64a282ffb1d2b9bb36188d0c"
$ws.Range("M25").Value = "Bernard"
$ws.Range("N25").Value = 4209403
$ws.Range("O25").Value = "'107491"
$ws.Range("P25").Value = "Approved"

# --- Apply the existing date style (numFmtId 14, same as column E elsewhere) to the new date cells ---
$ws.Range("E17").Copy()
$ws.Range("E18:E25").PasteSpecial(-4122)

$excel.CutCopyMode = 0